$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.515.96'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.587.84'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.91%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').Value = '2.596.61'
$ws.Range('E9').Value = '  -2.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('E12').Value = '  +10.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.357'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.85%  '
$ws.Range('D14').Value = '3.040.44'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.93%  '
$ws.Range('D16').Value = '59.501.56'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '2.592.76'
$ws.Range('E18').Value = '  -2.84%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '339.06'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.472'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.22'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.09'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.869'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '295.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.68'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '138.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0976'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.595'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.64'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('E48').Value = '  -3.04%  '
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.969.78'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.27%  '
